$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values that changed (rows 10-20) ---
$ws.Range("A10").Value = "CAO"
$ws.Range("B10").Value = "CAOA"

$ws.Range("A11").Value = "CHE"
$ws.Range("B11").Value = "CHERY"
$ws.Range("C11").Value = "CIRINA77"
$ws.Range("D11").Value = "Aline"

$ws.Range("A12").Value = "NIS"
$ws.Range("B12").Value = "NISSAN"
$ws.Range("C12").Value = "CIRINA77"
$ws.Range("D12").Value = "Aline"

$ws.Range("A13").Value = "ARG"
$ws.Range("B13").Value = "TENNECO ARGENTINA"
$ws.Range("C13").Value = "CIRINA77"
$ws.Range("D13").Value = "Aline"

$ws.Range("A14").Value = "VWC"
$ws.Range("B14").Value = "VWCO"

$ws.Range("A15").Value = "SCN"
$ws.Range("B15").Value = "SCANIA"

$ws.Range("A16").Value = "DAI"
$ws.Range("B16").Value = "MERCEDES BENZ"

$ws.Range("A17").Value = "CNH"
$ws.Range("B17").Value = "CNH"
$ws.Range("C17").Value = "ELTOBORG"

$ws.Range("A18").Value = "MWM"
$ws.Range("B18").Value = "MWM"

$ws.Range("A19").Value = "VOL"
$ws.Range("B19").Value = "VOLVO"

$ws.Range("A20").Value = "IVE"
$ws.Range("B20").Value = "IVECO"
$ws.Range("C20").Value = "ELTOBORG"
$ws.Range("D20").Value = "Elton"

# --- Fix cell formatting (fill colors) to match new grouping ---
# Rows 11-13 move from the "orange group" style to the "gold group" style (same as rows 2-10)
$ws.Range("A2:D2").Copy()
$ws.Range("A11:D13").PasteSpecial(-4122)

# Row 20 moves from the unique "Erik" style to the shared "Elton" style (same as rows 14-19)
$ws.Range("A14:D14").Copy()
$ws.Range("A20:D20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the selected cell/range shown in the sheet view ---
$ws.Range("B5").Select()
